# Applies the "Finished lit review" edit to Paper Notes.docx:
#  1. Splits the "But, these are..." bullet into two runs, wrapping
#     "But," in proofErr gramStart/gramEnd (as Word's grammar checker
#     would mark it).
#  2. Splits the "Maybe not a good paper..." bullet into three runs,
#     wrapping "movement based" in proofErr gramStart/gramEnd, while
#     preserving the _GoBack bookmark (it is relocated, see step 4).
#  3. Inserts the new "Data Preparation for Injury Prediction" section
#     (heading, bullet, blank italic paragraph, and the
#     "Effective Injury Prediction..." heading with a page-break marker).
#  4. Appends a final (now-empty) bulleted paragraph that carries the
#     relocated _GoBack bookmark, matching the end of the new content.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- locate the two target paragraphs by their current text -------------
$target1 = "But, these are on veteran football players (age >32), not high school and college students "
$target2 = "Maybe not a good paper for a direct comparison, but good to bring attention to if we end up using those movement based variables in the study "

$p1 = $null
$p2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq $target1) { $p1 = $p }
    if ($t -eq $target2) { $p2 = $p }
}

if ($p1 -eq $null) { throw "paragraph 1 not found" }
if ($p2 -eq $null) { throw "paragraph 2 not found" }

# --- step 1: split "But, these are..." into two runs with proofErr ------
$r1 = $p1.Range
$frag1 = $pkgOpen + '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>But,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> these are on veteran football players (age &gt;32), not high school and college students </w:t></w:r></w:p></w:body>' + $pkgClose
$r1.InsertXML($frag1)

# --- step 2: split "Maybe not a good paper..." into three runs ----------
# (the trailing _GoBack bookmark inside this paragraph is intentionally
#  dropped here and re-created on the new trailing paragraph in step 4)
$r2 = $p2.Range
$frag2 = $pkgOpen + '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Maybe not a good paper for a direct comparison, but good to bring attention to if we end up using those </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>movement based</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> variables in the study </w:t></w:r></w:p></w:body>' + $pkgClose
$r2.InsertXML($frag2)

# --- steps 3 & 4: insert the new section + trailing bookmark paragraph --
$insertPos = $p2.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$newBlock = $pkgOpen + '<w:body>' + `
  '<w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Data Preparation for Injury Prediction</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The stuff in here is pretty obvious, but it' + [char]0x2019 + 's another source we can point to </w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr></w:p>' + `
  '<w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Effective Injury Prediction in Professional Soccer with GPS Data and Machine Learning</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  '</w:body>' + $pkgClose

$insertRange.InsertXML($newBlock)
